# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (AB1) onto the new
# header cells so they keep the same bold/centered/bordered header style.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# New header labels
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every player row (2-42) gets the same team season record values.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 77   # AC: Wins
    $ws.Cells.Item($r, 30).Value = 84   # AD: Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE: Ties
}
